$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1574.1522668947819
$ws.Range("C1").Value = 1492.663573085847
$ws.Range("D1").Value = 1763.0141071223763
$ws.Range("E1").Value = 1504.6731898238745
$ws.Range("F1").Value = 1622.1153483992466
$ws.Range("G1").Value = 1551.0315107734464
$ws.Range("B2").Value = 1708.7513136991324
$ws.Range("C2").Value = 1557.8770301624133
$ws.Range("D2").Value = 1766.3218258974653
$ws.Range("E2").Value = 1542.9275929549901
$ws.Range("F2").Value = 1622.1153483992466
$ws.Range("G2").Value = 1548.5863743928262
$ws.Range("B3").Value = 1726.8072833923989
$ws.Range("C3").Value = 1602.9628770301626
$ws.Range("D3").Value = 1761.360247734832
$ws.Range("E3").Value = 1556.4760273972602
$ws.Range("F3").Value = 1655.5850988700565
$ws.Range("G3").Value = 1612.1599202889527
$ws.Range("B4").Value = 1788.3617255285349
$ws.Range("C4").Value = 1641.6078886310906
$ws.Range("D4").Value = 1873.8226860878542
$ws.Range("E4").Value = 1626.6090998043051
$ws.Range("F4").Value = 1716.2976694915253
$ws.Range("G4").Value = 1668.3980570432184
$ws.Range("B5").Value = 1855.6612489307101
$ws.Range("C5").Value = 1713.2621809744783
$ws.Range("D5").Value = 1930.0539052643651
$ws.Range("E5").Value = 1699.9300391389431
$ws.Range("F5").Value = 1740.4270244821091
$ws.Range("G5").Value = 1727.0813301781045
$ws.Range("B6").Value = 1909.8291580105097
$ws.Range("C6").Value = 1822.7563805104412
$ws.Range("D6").Value = 2037.554765454754
$ws.Range("E6").Value = 1771.6570450097845
$ws.Range("F6").Value = 1839.2795433145009
$ws.Range("G6").Value = 1814.2911944202269
$ws.Range("B7").Value = 1951.686178663082
$ws.Range("C7").Value = 1886.3596287703019
$ws.Range("D7").Value = 2142.574836563826
$ws.Range("E7").Value = 1844.1810176125243
$ws.Range("F7").Value = 1891.4300847457628
$ws.Range("G7").Value = 1898.2408768215221
$ws.Range("B8").Value = 2063.304900403275
$ws.Range("C8").Value = 1943.5220417633413
$ws.Range("D8").Value = 2195.4983369652482
$ws.Range("E8").Value = 1900.7656555772992
$ws.Range("F8").Value = 1970.8234463276835
$ws.Range("G8").Value = 1919.4320587868976
$ws.Range("B9").Value = 2130.6044238054506
$ws.Range("C9").Value = 2061.067285382831
$ws.Range("D9").Value = 2249.2487670604428
$ws.Range("E9").Value = 2004.3713307240703
$ws.Range("F9").Value = 2042.4331450094162
$ws.Range("G9").Value = 2035.983559596463
$ws.Range("B10").Value = 2248.7889527068314
$ws.Range("C10").Value = 2189.0788863109051
$ws.Range("D10").Value = 2377.4228695951369
$ws.Range("E10").Value = 2129.4951076320935
$ws.Range("F10").Value = 2181.760710922787
$ws.Range("G10").Value = 2127.2686511396191
$ws.Range("B11").Value = 2443.3009898570208
$ws.Range("C11").Value = 2343.6589327146175
$ws.Range("D11").Value = 2547.7703865122148
$ws.Range("E11").Value = 2287.294520547945
$ws.Range("F11").Value = 2398.1465395480227
$ws.Range("G11").Value = 2337.5503798729605
$ws.Range("B12").Value = 2615.6534278382014
$ws.Range("C12").Value = 2508.7053364269145
$ws.Range("D12").Value = 2708.1947471040257
$ws.Range("E12").Value = 2426.7636986301368
$ws.Range("F12").Value = 2526.5769774011301
$ws.Range("G12").Value = 2463.0673807448002
$ws.Range("B13").Value = 2784.7229622387877
$ws.Range("C13").Value = 2667.31090487239
$ws.Range("D13").Value = 2924.850326872348
$ws.Range("E13").Value = 2600.502446183953
$ws.Range("F13").Value = 2683.8069679849341
$ws.Range("G13").Value = 2626.0764727861506
$ws.Range("B14").Value = 2946.4059635830381
$ws.Range("C14").Value = 2828.331786542924
$ws.Range("D14").Value = 3094.370914095653
$ws.Range("E14").Value = 2755.1139921722111
$ws.Range("F14").Value = 2845.7071563088512
$ws.Range("G14").Value = 2768.7094283223319
$ws.Range("B15").Value = 3120.3998533545155
$ws.Range("C15").Value = 2964.394431554525
$ws.Range("D15").Value = 3260.5837825438698
$ws.Range("E15").Value = 2914.507338551859
$ws.Range("F15").Value = 2995.9318502824858
$ws.Range("G15").Value = 2928.4583385228548
$ws.Range("B16").Value = 3114.6547720884764
$ws.Range("C16").Value = 3050.5406032482601
$ws.Range("D16").Value = 3388.757885078564
$ws.Range("E16").Value = 2972.6859099804301
$ws.Range("F16").Value = 3094.0060028248586
$ws.Range("G16").Value = 3025.4487482874583
$ws.Range("B17").Value = 3275.517047537578
$ws.Range("C17").Value = 3167.2807424593971
$ws.Range("D17").Value = 3467.3162059869251
$ws.Range("E17").Value = 3088.2460861056747
$ws.Range("F17").Value = 3153.1618408662898
$ws.Range("G17").Value = 3130.5896126541293
$ws.Range("B18").Value = 3282.9035805939143
$ws.Range("C18").Value = 3193.8491879350354
$ws.Range("D18").Value = 3489.6433077187753
$ws.Range("E18").Value = 3152.8003913894322
$ws.Range("F18").Value = 3250.4576271186438
$ws.Range("G18").Value = 3178.6772948063272
$ws.Range("B19").Value = 3336.2507637785657
$ws.Range("C19").Value = 3301.7331786542927
$ws.Range("D19").Value = 3541.7398784264251
$ws.Range("E19").Value = 3159.9730919765161
$ws.Range("F19").Value = 3287.8192090395478
$ws.Range("G19").Value = 3200.6835222319096
$ws.Range("B20").Value = 3322.2984235610411
$ws.Range("C20").Value = 3299.3178654292346
$ws.Range("D20").Value = 3574.8170661773138
$ws.Range("E20").Value = 3187.8669275929547
$ws.Range("F20").Value = 3295.6028719397364
$ws.Range("G20").Value = 3251.2163407647281
$ws.Range("B21").Value = 3333.7885860931201
$ws.Range("C21").Value = 3312.1995359628777
$ws.Range("D21").Value = 3585.5671521963527
$ws.Range("E21").Value = 3218.1516634050877
$ws.Range("F21").Value = 3311.9485640301318
$ws.Range("G21").Value = 3256.9216589861753
$ws.Range("B22").Value = 3366.617621899059
$ws.Range("C22").Value = 3319.4454756380514
$ws.Range("D22").Value = 3571.5093474022251
$ws.Range("E22").Value = 3210.1819960861053
$ws.Range("F22").Value = 3315.8403954802261
$ws.Range("G22").Value = 3268.3322954290697
$ws.Range("B23").Value = 3360.87254063302
$ws.Range("C23").Value = 3304.1484918793508
$ws.Range("D23").Value = 3593.8364491340749
$ws.Range("E23").Value = 3242.0606653620348
$ws.Range("F23").Value = 3309.6134651600755
$ws.Range("G23").Value = 3270.7774318096904

$ws.Range("B1:G23").Select()
